# Clean up shared-string cell values that contain embedded line breaks,
# collapsing them to single-line text with a plain space in place of each
# newline. This matches the author's cleanup pass across the two influenza
# vaccine sheets ("Pediatric Influenza Vaccine" and "Adult Influenza Vaccine").

$wb = $excel.ActiveWorkbook

$ws3 = $wb.Worksheets.Item("Pediatric Influenza Vaccine ")
$ws3.Range("B3").Value  = "Fluzone Pediatric dose No Preservative"
$ws3.Range("B6").Value  = "Fluarix Preservative-Free"
$ws3.Range("B9").Value  = "FluMist No Preservative"
$ws3.Range("B10").Value = "Afluria No Preservative"
$ws3.Range("H10").Value = "Merck (CSL product)"

$ws4 = $wb.Worksheets.Item("Adult Influenza Vaccine ")
$ws4.Range("B5").Value  = "Agriflu No Preservative"
$ws4.Range("B7").Value  = "Fluvirin Preservative-free"
$ws4.Range("B8").Value  = "Fluarix Preservative-free"
$ws4.Range("B10").Value = "Flumist No Preservative"
